$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Window view size (cosmetic, matches author's saved window state)
$excel.Width = 1087

# Tweak two existing Tips cells
$ws.Range("D322").Value = "数组实现：待补充  defaultdict字典实现：总体来看是个字典套字典的形式eg： {a:{b:{c:{d}}}}，写一个自己的TreeNode类，属性有children和isWord，isWord在插入方法结束时候设置成True， 搜索方法就是遍历字符串 一个个字母去字典里搜，如果有一个字母不存在就返回False，因为都是唯一的一一对应，遍历字符串以最后一个字符的isWord来判断是否找到，startsWith方法则类似，只不过不需要通过isWord来判断。"
$ws.Range("D343").Value = "cmp_to_key：重写sort 第一步先把 nums 中的所有数字转字符串，形成字符串数组 nums_str 第二步比较两个字符串 x,y 拼接的结果 x + y 和 y + x 哪个更大，从而确定 x 和 y 谁排在前面然后将 nums_str 降序排序，第三步把整个数组排序的结果拼接成一个字符串，并返回  快速排序：字符串拼接比较大小 （a+b > b+a） + 快排倒序 # corner case: [0,0] return -> ""0"""

# Fill in category for row 499
$ws.Range("E499").Value = "bitManipulation"

# New daily entries rows 510-527
# row 510
$ws.Range("A510").Value = "872. Leaf-Similar Trees"
$ws.Range("B510").Value = "easy"
$ws.Range("C510").Value = "2021-11-20"
$ws.Range("D510").Value = "Tree(preorder)：遍历两棵树当到达叶子节点的时候加入res数组 然后对比两个res数组是否完全相同"
$ws.Range("E510").Value = "Tree(preorder)"
$ws.Rows(510).RowHeight = 75.1

# row 511
$ws.Range("A511").Value = "8. String to Integer (atoi)"
$ws.Range("B511").Value = "medium"
$ws.Range("C511").Value = "2021-11-20"
$ws.Range("D511").Value = "模拟：第一步check leading whitespace 第二步check neg or pos 第三步string to integer（字符可以直接进行比较 eg：s[j] < '0' or s[j] > '9'表示字符不是数字） 第四步范围检查  正则表达式：[ ]* 任意空白 ([+-]?\d+) []匹配一个字符 ？匹配0或1个正负号 至少一个字符的任意0-9的数字 ( )表示group                                    #https://leetcode-cn.com/problems/string-to-integer-atoi/solution/xiong-mao-shua-ti-python3-yi-qi-xue-xi-zheng-ze-bi/  https://www.liaoxuefeng.com/wiki/1016959663602400/1017639890281664"
$ws.Range("E511").Value = "正则表达式"
$ws.Rows(511).RowHeight = 100

# row 512
$ws.Range("A512").Value = "559. Maximum Depth of N-ary Tree"
$ws.Range("B512").Value = "easy"
$ws.Range("C512").Value = "2021-11-21"
$ws.Range("D512").Value = "DFS：写法一 维护全局变量 / 写法二  维护局部变量ans 寻找最深的child的深度 然后返回给上一层计算 BFS：层序遍历维护step"
$ws.Range("E512").Value = "Tree"
$ws.Range("F512").Value = 104
$ws.Rows(512).RowHeight = 75.1

# row 513
$ws.Range("A513").Value = "208. Implement Trie (Prefix Tree)"
$ws.Range("B513").Value = "meidum"
$ws.Range("C513").Value = "2021-11-21"
$ws.Range("D513").Value = "dict字典实现： insert：如果key不存在 就cur[w] = {} 然后移动cur 如果存在就直接移动cur 到末尾了以后设置一个cur['#'] = {} search：如果key不存在就返回false 如果存在一直移动cur 最后以'#'是否在cur中来判断是否找到 startsWith：如果key不存在就返回false 最后返回true"
$ws.Range("E513").Value = "Trie/数据结构"
$ws.Rows(513).RowHeight = 75.1

# row 514
$ws.Range("A514").Value = "677. Map Sum Pairs"
$ws.Range("B514").Value = "medium"
$ws.Range("C514").Value = "2021-11-21"
$ws.Range("D514").Value = "bruteForce：直接存储键值对然后sum的时候再扫描一个个key是否匹配来计算 Trie：类似与构建字典树 不过在insert的最后结尾的时候用cur['#']=val来表示值 要sum的时候先确定前缀是存在的 然后再把同一前缀下的所有val累加起来 写法一 DFS 具体做法就是用dfs 维护一个局部变量ans 然后先查看cur中key是否有""#"" 如果有就ans += cur[k] 如果没有就继续往下递归看它的val是否还有包含""#""的  / 写法二 BFS 用队列装下cur 然后一个个弹出来查看 与dfs类似只不过用了队列来操作"
$ws.Range("E514").Value = "Trie/数据结构"
$ws.Rows(514).RowHeight = 86

# row 515
$ws.Range("A515").Value = "剑指 Offer 45. 把数组排成最小的数"
$ws.Range("B515").Value = "medium"
$ws.Range("C515").Value = "2021-11-21"
$ws.Range("D515").Value = "cmp_to_key：重写sort 第一步先把 nums 中的所有数字转字符串，形成字符串数组 nums_str 第二步比较两个字符串 x,y 拼接的结果 x + y 和 y + x 哪个更小，从而确定 x 和 y 谁排在前面 ，具体实现用cmp_to_key方法(return 1表示大于) 第三步把整个数组排序的结果拼接成一个字符串，并返回   快速排序：改写快排 在比较的时候 str(nums[j]) + str(nums[low]) < str(nums[low]) + str(nums[j])则表示与正常快排中nums[j]<nums[low]一个意思 所以需要交换 #┭┮﹏┭┮不会证明 待学习   #cmp_to_key自定义函数需要:1.接收两个参数 p1 , p2; 2.返回1、0或-1，其中1代表 p1 > p2, 0代表 p1 == p2, -1代表 p1 < p2  #corner case: [0, 0] return -> ""00"""
$ws.Range("E515").Value = "sort"
$ws.Range("F515").Value = 179
$ws.Rows(515).RowHeight = 116

# row 516
$ws.Range("A516").Value = "2078. Two Furthest Houses With Different Colors"
$ws.Range("B516").Value = "easy"
$ws.Range("C516").Value = "2021-11-21"
$ws.Range("D516").Value = "bruteForce：暴力枚举维护maxLen"
$ws.Range("E516").Value = "枚举"
$ws.Rows(516).RowHeight = 75.1

# row 517
$ws.Range("A517").Value = "2079. Watering Plants"
$ws.Range("B517").Value = "medium"
$ws.Range("C517").Value = "2021-11-21"
$ws.Range("D517").Value = "模拟：题目给出max(plants[i]) <= capacity 那么就不需要考虑需要多次浇水的情况  定义total为目前桶内的水量 顺序遍历plants 如果plants[i] > total说明水不够用了 那么就回到河边并回来 因为是预判下一个 所以想补充水时的位置是前一个位置 即所用的步数为step += i + i（来回） 并且total = capacity（装满水）如果水够用就total-plants[i] 然后step+=1 最后返回step即可"
$ws.Range("E517").Value = "imitation/simulation"
$ws.Rows(517).RowHeight = 75.1

# row 518
$ws.Range("A518").Value = "2080. Range Frequency Queries"
$ws.Range("B518").Value = "medium"
$ws.Range("C518").Value = "2021-11-21"
$ws.Range("D518").Value = "二分：用二分插入获得idx然后相减就可以得到答案 #不一定每道题都是二分左右界和二分答案 二分插入可以用于求一个idx差值"
$ws.Range("E518").Value = "binarySearch"
$ws.Rows(518).RowHeight = 75.1

# row 519
$ws.Range("A519").Value = "384. Shuffle an Array"
$ws.Range("B519").Value = "medium"
$ws.Range("C519").Value = "2021-11-22"
$ws.Range("D519").Value = "概率：根据 1/n = (n - 1)/n * 1/(n - 1) = ... 实现等概率算法  做一个tmp为nums的copy 用for循环依次选择卡片 然后在所有卡片中随机选一张卡牌(1/n)与其进行调换 每次调换完 缩小选取随机卡片的范围继续调换剩余的卡片((n - 1)/n * 1/(n - 1)) 直到所有的都被调换过一次"
$ws.Range("E519").Value = "Math"
$ws.Rows(519).RowHeight = 75.1

# row 520
$ws.Range("A499:E499").Copy()
$ws.Range("A520:E520").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A520").Value = "54. Spiral Matrix"
$ws.Range("B520").Value = "medium"
$ws.Range("C520").Value = "2021-11-22"
$ws.Range("D520").Value = "DFS：设计函数为dfs(x, y, directionIdx, step) 进入函数就表示合法 那么意味着我们需要将matrix[x][y]加入res并设置为已经visited 再进行下一步前先用一个isValid函数验证下 下一步是否合理 如果合理就继续走不合理就 directionIdx=(directionIdx+1)%4 然后将新的参数传入dfs继续递归 注意设置visited防止重复访问造成死循环 step==row*col就可以退出循环了 corner case：[1]        4重for循环：先定义top bottom left right四个边 然后依次遍历直到top > bottom or left > right or right < left or down < top 四种情况下break掉即可            遍历的具体为left->right+1 top+=1 top->down+1 right-=1 ..."
$ws.Range("E520").Value = "DFS"
$ws.Rows(520).RowHeight = 114

# row 521
$ws.Range("A521").Value = "2. Add Two Numbers"
$ws.Range("B521").Value = "medium"
$ws.Range("C521").Value = "2021-11-22"
$ws.Range("D521").Value = "模拟：题目难度不大 用cur1和cur2指向l1和l2 当cur1 or cur2的时候保持循环 需要注意几个点 第一就是两个数字长度不一样的情况下需要对cur分别判断 if cur1:cur1 = cur1.next如此 否则会nullPointer  第二就是如果遍历完发现carry还有值 最后需要单独处理下carry"
$ws.Range("E521").Value = "LinkedList"
$ws.Rows(521).RowHeight = 75.1

# row 522
$ws.Range("A522").Value = "322. Coin changes "
$ws.Range("B522").Value = "medium"
$ws.Range("C522").Value = "2021-11-22"
$ws.Range("D522").Value = "BFS：最短路问题反应出想到BFS 遍历当前层for coin in coins 一个个试过去  计算 left_amount 如果left_amount==0就返回步数说明找到了 这里可以做一个剪枝提前对coins排序然后left_amount<0就break掉 最后将left_amount加入set防止重复访问 队列里之推入没有访问过的left_amount 不重复推入1是因为后面的步数肯定更多没必要要找到前面就找到了2是因为会死循环 如果全遍历完都没找到就返回-1"
$ws.Range("E522").Value = "DFS(step)/BFS"
$ws.Rows(522).RowHeight = 75.1

# row 523
$ws.Range("A523").Value = "295. Find Median from Data Stream"
$ws.Range("B523").Value = "hard"
$ws.Range("C523").Value = "2021-11-23"
$ws.Range("D523").Value = "双heap：大根堆+小根堆 调用addNum的时候提前先判断左右堆内的元素数量 分为两种情况 第一种是偶数情况 此时元素数量相同 当lenRight == 0或者num小于小根堆的堆顶元素就直接加入大根堆 如果num>self.right[0]说明需要放到小根堆 此时需要先弹出小根堆的堆顶元素 将其加入大根堆 然后再把num加入小根堆   第二种是奇数情况 此时L比R多一个元素(这里定义要么大根堆和小根堆中元素数量相同要么大根堆里多一个只有这两种情况) 如果此时num大于等于大根堆的堆顶元素  那么就直接加入大根堆 如果小于大根堆的堆顶元素那么 就需要从大根堆中弹出堆顶元素 将其加入小根堆 然后把num加入大根堆 #注意tuple中的取值 在findMedian函数中如果大根堆和小根堆内元素数量相同就直接 float(两个堆的堆顶元素相加)/2 如果不等那么就返回大根堆的堆顶元素"
$ws.Range("E523").Value = "Heap"
$ws.Rows(523).RowHeight = 125

# row 524
$ws.Range("A524").Value = "859. Buddy Strings"
$ws.Range("A524").Style = $ws.Range("C524").Style
$ws.Range("B524").Value = "easy"
$ws.Range("C524").Value = "2021-11-23"
$ws.Range("D524").Value = "哈希表：1.当 s 与 goal 长度 或 词频不同，必然不为亲密字符；2.当「s 与 goal 不同的字符数量为 2（能够相互交换）」或「s 与 goal 不同的字符数量为 0，但同时 s 中有出现数量超过 2 的字符（能够自己相互交换）」时，两者必然为亲密字符。 具体实现步骤为先排除长度不同的情况 注意还需要排除词频不同的情况  计算位置不同的字母个数为多少 在末尾验证diff == 2  如果不是2个就要做另外一个判断 验证s中有出现数量超过 2 的字符  eg1:s=""aa"" goal=""aa""  总来说得满足dif==2 or (diff==0 or s中有重复的字符)"
$ws.Range("E524").Value = "hashTable"
$ws.Rows(524).RowHeight = 100

# row 525
$ws.Range("A525").Value = "面试题 17.17. Multi Search LCCI"
$ws.Range("B525").Value = "medium"
$ws.Range("C525").Value = "2021-11-23"
$ws.Range("D525").Value = "Trie：这里反转思路我们将smalls插入字典树 末尾记录字符串，方便后面遍历 search的时候传入bigs[i:] 每次都判断cur中是否包含了'#' 如果包含就将这个string加入res 然后返回来res里面装着所有符合开头的字符串 然后用dic把他们分别存起来{string:[idx1, idx2...]} 这样最后再按照smalls的顺序遍历一次 用dic[small]就可以拿到顺序正确的答案了"
$ws.Range("E525").Value = "Trie"
$ws.Rows(525).RowHeight = 75.1

# row 526
$ws.Range("A526").Value = "232. Implement Queue using Stacks"
$ws.Range("B526").Value = "easy"
$ws.Range("C526").Value = "2021-11-24"
$ws.Range("D526").Value = "双栈：在pop前先判断 outStack是否为空 如果outStack为空就加入inStack的所有元素 如果不为空就直接从outStack中pop出元素"
$ws.Range("E526").Value = "Stack"
$ws.Rows(526).RowHeight = 75.1

# row 527
$ws.Range("A527").Value = "547. Number of Provinces"
$ws.Range("B527").Value = "medium"
$ws.Range("C527").Value = "2021-11-25"
$ws.Rows(527).RowHeight = 75.1

# Scroll/selection state at save time
$ws.Range("D527").Select()
$excel.ActiveWindow.ScrollRow = 524
$excel.ActiveWindow.FreezePanes = $true
